$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Column widths (B: 33 -> 41, H: 12 -> 18) ---------------------------
# ColumnWidth uses "character" units; the raw OOXML <col width> value
# Excel stores is ColumnWidth + 5/6 (default Calibri 11 padding), so we
# back that offset out to land on an exact integer in the saved file.
$ws.Columns.Item(2).ColumnWidth = 41 - 5/6
$ws.Columns.Item(8).ColumnWidth = 18 - 5/6

# --- Drop every existing hyperlink so we can rebuild the list cleanly --
# (Hyperlinks.Add on a cell that already carries a hyperlink duplicates
# the relationship instead of replacing it, so the list is rebuilt from
# scratch below.)
$ws.Range("A1").Hyperlinks.Delete()

# --- Row 2: existing entry, only the capture timestamp changes --------
$ws.Range("A2").Value = "2026-01-16 01:27:53"

# --- Row 3: now the Python/Docker scraping listing ---------------------
$ws.Range("A3").Value = "2026-01-16 01:27:53"
$ws.Range("B3").Value = "【募集】Python / Docker 日次データ スクレイピングシステム構築"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5469627"
$ws.Range("G3").Value = 248
$ws.Range("H3").Value = "🔥Python ◆スクレイピング"

# --- Row 4 (new): RoboTANGO RPA listing ---------------------------------
$ws.Range("A4").Value = "2026-01-16 01:27:53"
$ws.Range("B4").Value = "【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5405023"
$ws.Range("G4").Value = 178
$ws.Range("H4").Value = "★bot ◆ツール"

# --- Row 5 (new): the GAS project listing, pushed down from old row 3 --
$ws.Range("A5").Value = "2026-01-16 01:27:53"
$ws.Range("B5").Value = "【急募】GASプロジェクトの作成依頼"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5471552"
$ws.Range("G5").Value = 25

# --- Rebuild hyperlinks for every URL cell in order ---------------------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5471035")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5469627")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5405023")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5471552")

# Hyperlinks.Add stamps its own (duplicate-looking) style; re-apply the
# shared "Hyperlink" cell style afterwards so every link cell points at
# the same style record the workbook already had for F2/F3.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
